# Auto-generated edit script: updates currentAveragePrice-derived columns
# (H..N) across multiple sheets to reflect refreshed market data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 960.3
$ws.Range("I9").Value = 800.4286
$ws.Range("K9").Value = 800.4286
$ws.Range("M9").Value = -631.4286
$ws.Range("H26").Value = 3756.25
$ws.Range("I26").Value = 1675
$ws.Range("J26").Value = 10000
$ws.Range("K26").Value = 1675
$ws.Range("L26").Value = 10000
$ws.Range("M26").Value = -1331
$ws.Range("N26").Value = -10688
$ws.Range("H74").Value = 7880.8125
$ws.Range("I74").Value = 7672.8667
$ws.Range("K74").Value = 7672.8667
$ws.Range("M74").Value = -6736.8667
$ws.Range("H77").Value = 7880.8125
$ws.Range("I77").Value = 7672.8667
$ws.Range("K77").Value = 38364.3335
$ws.Range("M77").Value = -33684.3335
$ws.Range("H103").Value = 789.6
$ws.Range("I103").Value = 724
$ws.Range("J103").Value = 833.3333
$ws.Range("K103").Value = 2172
$ws.Range("L103").Value = 2499.9999
$ws.Range("M103").Value = -1586
$ws.Range("N103").Value = -3671.9999
$ws.Range("H106").Value = 15354.637
$ws.Range("I106").Value = 13990.3
$ws.Range("J106").Value = 28998
$ws.Range("K106").Value = 13990.3
$ws.Range("L106").Value = 28998
$ws.Range("M106").Value = -13359.3
$ws.Range("N106").Value = -30260
$ws.Range("H112").Value = 3661.2856
$ws.Range("I112").Value = 945
$ws.Range("J112").Value = 4747.8
$ws.Range("K112").Value = 2835
$ws.Range("L112").Value = 14243.4
$ws.Range("M112").Value = -1727
$ws.Range("N112").Value = -16459.4
$ws.Range("H113").Value = 7877.4
$ws.Range("I113").Value = 6943.5
$ws.Range("J113").Value = 8500
$ws.Range("K113").Value = 6943.5
$ws.Range("L113").Value = 8500
$ws.Range("M113").Value = -3689.5
$ws.Range("N113").Value = -15008
$ws.Range("H132").Value = 83340750
$ws.Range("I132").Value = 83340750
$ws.Range("K132").Value = 250022250
$ws.Range("M132").Value = -250019720
$ws.Range("H135").Value = 764.9167
$ws.Range("I135").Value = 792.63635
$ws.Range("K135").Value = 7133.72715
$ws.Range("M135").Value = -4598.72715
$ws.Range("H138").Value = 3941.6875
$ws.Range("I138").Value = 4783.5835
$ws.Range("K138").Value = 14350.7505
$ws.Range("M138").Value = -9210.750499999998
$ws.Range("H141").Value = 3282
$ws.Range("I141").Value = 2993
$ws.Range("K141").Value = 8979
$ws.Range("M141").Value = -3799

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6219.5
$ws.Range("I61").Value = 7326.3335
$ws.Range("K61").Value = 7326.3335
$ws.Range("M61").Value = -7114.3335
$ws.Range("H74").Value = 4888
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 4888
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 4888
$ws.Range("M74").Value = ""
$ws.Range("N74").Value = -6636
$ws.Range("H77").Value = 4888
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 4888
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 24440
$ws.Range("M77").Value = ""
$ws.Range("N77").Value = -33176
$ws.Range("H102").Value = 2882.7334
$ws.Range("I102").Value = 2882.7334
$ws.Range("K102").Value = 2882.7334
$ws.Range("M102").Value = -1260.7334
$ws.Range("H132").Value = 4171.0415
$ws.Range("I132").Value = 4238.826
$ws.Range("K132").Value = 12716.478
$ws.Range("M132").Value = -10186.478
$ws.Range("H136").Value = 6219.5
$ws.Range("I136").Value = 7326.3335
$ws.Range("K136").Value = 21979.0005
$ws.Range("M136").Value = -19429.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 612
$ws.Range("I22").Value = 634.1818
$ws.Range("J22").Value = 563.2
$ws.Range("K22").Value = 634.1818
$ws.Range("L22").Value = 563.2
$ws.Range("M22").Value = -461.1818
$ws.Range("N22").Value = -909.2
$ws.Range("H86").Value = 22458.223
$ws.Range("I86").Value = 22343.818
$ws.Range("K86").Value = 22343.818
$ws.Range("M86").Value = -21220.818
$ws.Range("H89").Value = 22458.223
$ws.Range("I89").Value = 22343.818
$ws.Range("K89").Value = 111719.09
$ws.Range("M89").Value = -106103.09
$ws.Range("H134").Value = 1030.4
$ws.Range("I134").Value = 1039.7142
$ws.Range("K134").Value = 3119.1426
$ws.Range("M134").Value = -584.1425999999997

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2180.0833
$ws.Range("I58").Value = 2306.25
$ws.Range("K58").Value = 2306.25
$ws.Range("M58").Value = -2103.25
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").Value = ""
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").Value = ""
$ws.Range("H99").Value = 3443
$ws.Range("J99").Value = 4008.6667
$ws.Range("L99").Value = 4008.6667
$ws.Range("N99").Value = -7004.6667
$ws.Range("H126").Value = 3443
$ws.Range("J126").Value = 4008.6667
$ws.Range("L126").Value = 12026.0001
$ws.Range("N126").Value = -16966.0001
$ws.Range("H132").Value = 12508198
$ws.Range("I132").Value = 15394419
$ws.Range("J132").Value = 1237
$ws.Range("K132").Value = 46183257
$ws.Range("L132").Value = 3711
$ws.Range("M132").Value = -46180727
$ws.Range("N132").Value = -8771
$ws.Range("H136").Value = 2180.0833
$ws.Range("I136").Value = 2306.25
$ws.Range("K136").Value = 6918.75
$ws.Range("M136").Value = -4368.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 2521.2144
$ws.Range("I113").Value = 2386.4285
$ws.Range("J113").Value = 2656
$ws.Range("K113").Value = 7159.2855
$ws.Range("L113").Value = 7968
$ws.Range("M113").Value = -4989.2855
$ws.Range("N113").Value = -12308

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 24497.5
$ws.Range("J95").Value = 38995
$ws.Range("L95").Value = 38995
$ws.Range("N95").Value = -44487
$ws.Range("H102").Value = 1600.4348
$ws.Range("I102").Value = 1648
$ws.Range("K102").Value = 1648
$ws.Range("M102").Value = -26
$ws.Range("H132").Value = 47623332
$ws.Range("I132").Value = 4995.5
$ws.Range("K132").Value = 14986.5
$ws.Range("M132").Value = -12456.5
$ws.Range("H136").Value = 82500
$ws.Range("J136").Value = 82500
$ws.Range("L136").Value = 247500
$ws.Range("N136").Value = -252600

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 47623300
$ws.Range("I40").Value = 83336130
$ws.Range("K40").Value = 83336130
$ws.Range("M40").Value = -83335994
$ws.Range("H48").Value = 19139.8
$ws.Range("I48").Value = 19139.8
$ws.Range("K48").Value = 19139.8
$ws.Range("M48").Value = -18478.8
$ws.Range("H122").Value = 2914.3333
$ws.Range("J122").Value = 3499.5
$ws.Range("L122").Value = 10498.5
$ws.Range("N122").Value = -15398.5
$ws.Range("H132").Value = 5499.6
$ws.Range("I132").Value = 5332.6665
$ws.Range("J132").Value = 5750
$ws.Range("K132").Value = 15997.9995
$ws.Range("L132").Value = 17250
$ws.Range("M132").Value = -13467.9995
$ws.Range("N132").Value = -22310

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 76940780
$ws.Range("I132").Value = 24627.777
$ws.Range("J132").Value = 250002140
$ws.Range("K132").Value = 73883.33099999999
$ws.Range("L132").Value = 750006420
$ws.Range("M132").Value = -71353.33099999999
$ws.Range("N132").Value = -750011480
